$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (item 4): cập nhật ngày bắt đầu/kết thúc và ngày thực tế ---
$ws.Range("E6").Value = 43375
$ws.Range("F6").Value = 43377
$ws.Range("G6").Value = 43375
$ws.Range("H6").Value = 43377

# G6/H6 previously had no date in them; give them the same date number
# format (style) already used by the rest of the date columns (copy from E6).
$ws.Range("E6").Copy()
$ws.Range("G6:H6").PasteSpecial(-4122)

# --- Row 7 (item 5) ---
$ws.Range("E7").Value = 43380
$ws.Range("F7").Value = 43381
$ws.Range("G7").Value = 43380
$ws.Range("H7").Value = 43388

$ws.Range("E7").Copy()
$ws.Range("G7:H7").PasteSpecial(-4122)

# --- Row 8 (item 6) ---
$ws.Range("G8").Value = 43390

$excel.CutCopyMode = 0

# --- Sheet view: top-left visible column C, selection on G9 ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G9").Select()
